{"js": "// Replace the 15 lattice-multiplication problems (5 rows x 3 cols) in the\n// document's single table with a new set of problems, keeping each cell's\n// existing layout/formatting (font size run, <w:br/> separated lines).\n//\n// Each cell encodes one problem \"A x B\" (two 2-digit numbers) as 5 lines:\n//   1) \"A x B\"\n//   2) \"  b1    b2\"      (digits of B, spaced)\n//   3) \"  ----\"\n//   4) \"a1|    |\"        (first digit of A)\n//   5) \"a2|    |\"        (second digit of A)\n\nconst newProblems = [\n  \"43 x 90\", \"72 x 43\", \"30 x 10\",\n  \"90 x 56\", \"23 x 85\", \"65 x 91\",\n  \"42 x 59\", \"73 x 56\", \"26 x 92\",\n  \"58 x 38\", \"33 x 47\", \"60 x 68\",\n  \"40 x 26\", \"89 x 58\", \"97 x 94\"\n];\n\nfunction cellLines(problem) {\n  const [a, b] = problem.split(\" x \");\n  return [\n    problem,\n    \"  \" + b[0] + \"    \" + b[1],\n    \"  ----\",\n    a[0] + \"|    |\",\n    a[1] + \"|    |\"\n  ];\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const para = cell.body.paragraphs.items[0];\n    const text = cellLines(newProblems[idx]).join(\"\\u000b\");\n    para.insertText(text, Word.InsertLocation.replace);\n\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication problems (5 rows x 3 cols) in the\n# document's single table with a new set of problems, keeping each cell's\n# existing layout/formatting (font size run, line-break separated lines).\n#\n# Each cell encodes one problem \"A x B\" (two 2-digit numbers) as 5 lines\n# separated by a manual line break (chr(11), same as <w:br/>):\n#   1) \"A x B\"\n#   2) \"  b1    b2\"      (digits of B, spaced)\n#   3) \"  ----\"\n#   4) \"a1|    |\"        (first digit of A)\n#   5) \"a2|    |\"        (second digit of A)\n\n$d = $word.ActiveDocument\n\n$newProblems = @(\n  \"43 x 90\", \"72 x 43\", \"30 x 10\",\n  \"90 x 56\", \"23 x 85\", \"65 x 91\",\n  \"42 x 59\", \"73 x 56\", \"26 x 92\",\n  \"58 x 38\", \"33 x 47\", \"60 x 68\",\n  \"40 x 26\", \"89 x 58\", \"97 x 94\"\n)\n\n$brk = [char]11\n\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $problem = $newProblems[$idx]\n    $parts = $problem.Split(\" x \")\n    $a = $parts[0]\n    $b = $parts[1]\n\n    $a1 = $a.Substring(0,1)\n    $a2 = $a.Substring(1,1)\n    $b1 = $b.Substring(0,1)\n    $b2 = $b.Substring(1,1)\n\n    # NOTE: use string interpolation (not the `+` operator) to join these\n    # pieces \u2014 strings that look numeric (e.g. \"9\", \"  9\") get coerced to\n    # numbers by `+`, corrupting the spacing.\n    $line1 = \"$problem\"\n    $line2 = \"  ${b1}    ${b2}\"\n    $line3 = \"  ----\"\n    $line4 = \"${a1}|    |\"\n    $line5 = \"${a2}|    |\"\n\n    $newText = \"${line1}${brk}${line2}${brk}${line3}${brk}${line4}${brk}${line5}\"\n\n    $cell = $tbl.Cell($r, $c)\n    $cell.Range.Text = $newText\n\n    $idx++\n  }\n}\n"}
